$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new row of data (row 43) following the existing pattern.
# Force column A to be treated as text first so the "MM/DD/YYYY" string isn't
# auto-converted into a date serial number, then clear the formatting back to
# the workbook default so the new cell matches the unstyled cells above it.
$ws.Range("A43").NumberFormat = "@"
$ws.Range("A43").Value = "10/14/2025"
$ws.Range("A43").ClearFormats()

$ws.Range("B43").Value = 0.1768496803993662
$ws.Range("C43").Value = 0.8231503196006338
